$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D:G stay text-formatted like the original inline-string cells
$ws.Range("D2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "283.67"
$ws.Range("E2").Value = "2.37%"
$ws.Range("F2").Value = "12-1-2023"
$ws.Range("G2").Value = "0"

$ws.Range("D3").Value = "28.57"
$ws.Range("E3").Value = "5.19%"
$ws.Range("F3").Value = "12-1-2023"
$ws.Range("G3").Value = "0"

$ws.Range("D4").Value = "4.898"
$ws.Range("E4").Value = "0.47%"
$ws.Range("F4").Value = "12-1-2023"
$ws.Range("G4").Value = "0"

$ws.Range("D5").Value = "0.06510"
$ws.Range("E5").Value = "1.42%"
$ws.Range("F5").Value = "12-1-2023"
$ws.Range("G5").Value = "0"

$ws.Range("D6").Value = "7.152"
$ws.Range("E6").Value = "2.77%"
$ws.Range("F6").Value = "12-1-2023"
$ws.Range("G6").Value = "0"

$ws.Range("D7").Value = "1.285"
$ws.Range("E7").Value = "3.45%"
$ws.Range("F7").Value = "12-1-2023"
$ws.Range("G7").Value = "0"

$ws.Range("D8").Value = "0.9195"
$ws.Range("E8").Value = "3.94%"
$ws.Range("F8").Value = "12-1-2023"
$ws.Range("G8").Value = "0"

$ws.Range("D9").Value = "0.1557"
$ws.Range("E9").Value = "2.81%"
$ws.Range("F9").Value = "12-1-2023"
$ws.Range("G9").Value = "0"

$ws.Range("D10").Value = "0.06290"
$ws.Range("E10").Value = "22.89%"
$ws.Range("F10").Value = "12-1-2023"
$ws.Range("G10").Value = "0"

$ws.Range("D11").Value = "0.07558"
$ws.Range("E11").Value = "0.56%"
$ws.Range("F11").Value = "12-1-2023"
$ws.Range("G11").Value = "0"

$ws.Range("D12").Value = "0.02930"
$ws.Range("E12").Value = "-1.30%"
$ws.Range("F12").Value = "12-1-2023"
$ws.Range("G12").Value = "0"

$ws.Range("D13").Value = "0.08969"
$ws.Range("E13").Value = "-0.41%"
$ws.Range("F13").Value = "12-1-2023"
$ws.Range("G13").Value = "0"

$ws.Range("D14").Value = "0.001589"
$ws.Range("E14").Value = "0.74%"
$ws.Range("F14").Value = "12-1-2023"
$ws.Range("G14").Value = "0"

$ws.Range("D15").Value = "0.0006421"
$ws.Range("E15").Value = "0.11%"
$ws.Range("F15").Value = "12-1-2023"
$ws.Range("G15").Value = "0"

$ws.Range("D16").Value = "0.005980"
$ws.Range("E16").Value = "3.27%"
$ws.Range("F16").Value = "12-1-2023"
$ws.Range("G16").Value = "0"

$ws.Range("D17").Value = "3.517"
$ws.Range("E17").Value = "1.51%"
$ws.Range("F17").Value = "12-1-2023"
$ws.Range("G17").Value = "0"

$ws.Range("D18").Value = "3.346"
$ws.Range("E18").Value = "0.89%"
$ws.Range("F18").Value = "12-1-2023"
$ws.Range("G18").Value = "0"

$ws.Range("D19").Value = "2.238"
$ws.Range("E19").Value = "-1.51%"
$ws.Range("F19").Value = "12-1-2023"
$ws.Range("G19").Value = "0"

$ws.Range("D20").Value = "0.3147"
$ws.Range("E20").Value = "0.36%"
$ws.Range("F20").Value = "12-1-2023"
$ws.Range("G20").Value = "0"

$ws.Range("D21").Value = "0.1351"
$ws.Range("E21").Value = "1.13%"
$ws.Range("F21").Value = "12-1-2023"
$ws.Range("G21").Value = "0"

$ws.Range("D22").Value = "3.988"
$ws.Range("E22").Value = "1.89%"
$ws.Range("F22").Value = "12-1-2023"
$ws.Range("G22").Value = "0"

$ws.Range("D23").Value = "0.1565"
$ws.Range("E23").Value = "13.39%"
$ws.Range("F23").Value = "12-1-2023"
$ws.Range("G23").Value = "0"

$ws.Range("D24").Value = "0.04391"
$ws.Range("E24").Value = "-0.77%"
$ws.Range("F24").Value = "12-1-2023"
$ws.Range("G24").Value = "0"

$ws.Range("D25").Value = "0.001173"
$ws.Range("E25").Value = "-0.20%"
$ws.Range("F25").Value = "12-1-2023"
$ws.Range("G25").Value = "0"

$ws.Range("D26").Value = "0.004323"
$ws.Range("E26").Value = "11.96%"
$ws.Range("F26").Value = "12-1-2023"
$ws.Range("G26").Value = "0"

$ws.Range("F27").Value = "12-1-2023"
$ws.Range("G27").Value = "0"

$ws.Range("D28").Value = "0.0001176"
$ws.Range("E28").Value = "-2.03%"
$ws.Range("F28").Value = "12-1-2023"
$ws.Range("G28").Value = "0"

$ws.Range("D29").Value = "0.0001644"
$ws.Range("E29").Value = "-15.15%"
$ws.Range("F29").Value = "12-1-2023"
$ws.Range("G29").Value = "0"

$ws.Range("F30").Value = "12-1-2023"
$ws.Range("G30").Value = "0"

$ws.Range("F31").Value = "12-1-2023"
$ws.Range("G31").Value = "0"

$ws.Range("F32").Value = "12-1-2023"
$ws.Range("G32").Value = "0"

$ws.Range("F33").Value = "12-1-2023"
$ws.Range("G33").Value = "0"

$ws.Range("F34").Value = "12-1-2023"
$ws.Range("G34").Value = "0"

$ws.Range("F35").Value = "12-1-2023"
$ws.Range("G35").Value = "0"

$ws.Range("F36").Value = "12-1-2023"
$ws.Range("G36").Value = "0"

$ws.Range("F37").Value = "12-1-2023"
$ws.Range("G37").Value = "0"

$ws.Range("F38").Value = "12-1-2023"
$ws.Range("G38").Value = "0"

$ws.Range("F39").Value = "12-1-2023"
$ws.Range("G39").Value = "0"

$ws.Range("D40").Value = "0.04125"
$ws.Range("E40").Value = "-0.71%"
$ws.Range("F40").Value = "12-1-2023"
$ws.Range("G40").Value = "0"

$ws.Range("D41").Value = "0.006708"
$ws.Range("E41").Value = "-1.17%"
$ws.Range("F41").Value = "12-1-2023"
$ws.Range("G41").Value = "0"

$ws.Range("D42").Value = "0.1391"
$ws.Range("E42").Value = "18.28%"
$ws.Range("F42").Value = "12-1-2023"
$ws.Range("G42").Value = "0"

$ws.Range("D43").Value = "0.002064"
$ws.Range("E43").Value = "-13.75%"
$ws.Range("F43").Value = "12-1-2023"
$ws.Range("G43").Value = "0"

$ws.Range("D44").Value = "0.01190"
$ws.Range("E44").Value = "6.03%"
$ws.Range("F44").Value = "12-1-2023"
$ws.Range("G44").Value = "0"

$ws.Range("D45").Value = "0.00005530"
$ws.Range("E45").Value = "6.29%"
$ws.Range("F45").Value = "12-1-2023"
$ws.Range("G45").Value = "0"

$ws.Range("D46").Value = "1.628"
$ws.Range("E46").Value = "9.52%"
$ws.Range("F46").Value = "12-1-2023"
$ws.Range("G46").Value = "0"

$ws.Range("D47").Value = "0.01843"
$ws.Range("E47").Value = "-8.97%"
$ws.Range("F47").Value = "12-1-2023"
$ws.Range("G47").Value = "0"

$ws.Range("F48").Value = "12-1-2023"
$ws.Range("G48").Value = "0"

$ws.Range("F49").Value = "12-1-2023"
$ws.Range("G49").Value = "0"

$ws.Range("F50").Value = "12-1-2023"
$ws.Range("G50").Value = "0"

$ws.Range("F51").Value = "12-1-2023"
$ws.Range("G51").Value = "0"
